$wb = $excel.ActiveWorkbook

# --- Update header labels to include units ("Added units to input files") ---

# Cluster1 sheet: cu_p_ch_max -> "cu_p_ch_max (kW)", cu_p_ds_max -> "cu_p_ds_max (kW)"
$wsCluster = $wb.Worksheets.Item("Cluster1")
$wsCluster.Range("B1").Value = "cu_p_ch_max (kW)"
$wsCluster.Range("C1").Value = "cu_p_ds_max (kW)"

# Capacity1 sheet: LB -> "LB (kW)", UB -> "UB (kW)"
$wsCapacity = $wb.Worksheets.Item("Capacity1")
$wsCapacity.Range("B1").Value = "LB (kW)"
$wsCapacity.Range("C1").Value = "UB (kW)"

# --- Restore selections / active sheet to match the saved view state ---

$wsCapacity.Activate() | Out-Null
$wsCapacity.Range("B1").Select() | Out-Null

$wsCluster.Activate() | Out-Null
$wsCluster.Range("C1").Select() | Out-Null
